$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2285.2273
$ws.Range("I15").Value = 2285.2273
$ws.Range("K15").Value = 6855.6819
$ws.Range("M15").Value = -6686.6819
$ws.Range("H51").Value = 6487.5
$ws.Range("I51").Value = 4500
$ws.Range("J51").Value = 7150
$ws.Range("K51").Value = 4500
$ws.Range("L51").Value = 7150
$ws.Range("M51").Value = -4016
$ws.Range("N51").Value = -8118
$ws.Range("H137").Value = 10206135
$ws.Range("I137").Value = 14707936
$ws.Range("K137").Value = 44123808
$ws.Range("M137").Value = -44121258
$ws.Range("H138").Value = 1380.1357
$ws.Range("I138").Value = 1188.7307
$ws.Range("J138").Value = 1723.3448
$ws.Range("K138").Value = 3566.1921
$ws.Range("L138").Value = 5170.0344
$ws.Range("M138").Value = 1573.8079
$ws.Range("N138").Value = -15450.0344
$ws.Range("H141").Value = 9370.723
$ws.Range("I141").Value = 2410.5
$ws.Range("J141").Value = 65052.5
$ws.Range("K141").Value = 7231.5
$ws.Range("L141").Value = 195157.5
$ws.Range("M141").Value = -2051.5
$ws.Range("N141").Value = -205517.5

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 3004.5715
$ws.Range("I46").Value = 3400
$ws.Range("J46").Value = 2846.4
$ws.Range("K46").Value = 3400
$ws.Range("L46").Value = 2846.4
$ws.Range("M46").Value = -3081
$ws.Range("N46").Value = -3484.4
$ws.Range("H61").Value = 2862.1667
$ws.Range("I61").Value = 1754.2858
$ws.Range("J61").Value = 3831.5625
$ws.Range("K61").Value = 1754.2858
$ws.Range("L61").Value = 3831.5625
$ws.Range("M61").Value = -1542.2858
$ws.Range("N61").Value = -4255.5625
$ws.Range("H74").Value = 27779558
$ws.Range("I74").Value = 38463412
$ws.Range("J74").Value = 1537.3334
$ws.Range("K74").Value = 38463412
$ws.Range("L74").Value = 1537.3334
$ws.Range("M74").Value = -38462538
$ws.Range("N74").Value = -3285.3334
$ws.Range("H77").Value = 27779558
$ws.Range("I77").Value = 38463412
$ws.Range("J77").Value = 1537.3334
$ws.Range("K77").Value = 192317060
$ws.Range("L77").Value = 7686.666999999999
$ws.Range("M77").Value = -192312692
$ws.Range("N77").Value = -16422.667
$ws.Range("H132").Value = 2555174
$ws.Range("I132").Value = 1754.9642
$ws.Range("J132").Value = 5959733
$ws.Range("K132").Value = 5264.892599999999
$ws.Range("L132").Value = 17879199
$ws.Range("M132").Value = -2734.892599999999
$ws.Range("N132").Value = -17884259
$ws.Range("H136").Value = 2862.1667
$ws.Range("I136").Value = 1754.2858
$ws.Range("J136").Value = 3831.5625
$ws.Range("K136").Value = 5262.857400000001
$ws.Range("L136").Value = 11494.6875
$ws.Range("M136").Value = -2712.857400000001
$ws.Range("N136").Value = -16594.6875

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 586
$ws.Range("I64").Value = 600.1429000000001
$ws.Range("J64").Value = 577.75
$ws.Range("K64").Value = 600.1429000000001
$ws.Range("L64").Value = 577.75
$ws.Range("M64").Value = -375.1429000000001
$ws.Range("N64").Value = -1027.75
$ws.Range("H67").Value = 586
$ws.Range("I67").Value = 600.1429000000001
$ws.Range("J67").Value = 577.75
$ws.Range("K67").Value = 600.1429000000001
$ws.Range("L67").Value = 577.75
$ws.Range("M67").Value = 179.8570999999999
$ws.Range("N67").Value = -2137.75
$ws.Range("H134").Value = 5772.2705
$ws.Range("I134").Value = 2905.5
$ws.Range("K134").Value = 8716.5
$ws.Range("M134").Value = -6181.5

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6412198
$ws.Range("I31").Value = 1372.6904
$ws.Range("J31").Value = 13891494
$ws.Range("K31").Value = 1372.6904
$ws.Range("L31").Value = 13891494
$ws.Range("M31").Value = -1077.6904
$ws.Range("N31").Value = -13892084
$ws.Range("H34").Value = 6412198
$ws.Range("I34").Value = 1372.6904
$ws.Range("J34").Value = 13891494
$ws.Range("K34").Value = 1372.6904
$ws.Range("L34").Value = 13891494
$ws.Range("M34").Value = -1170.6904
$ws.Range("N34").Value = -13891898
$ws.Range("H58").Value = 1569157.1
$ws.Range("I58").Value = 1358.8096
$ws.Range("J58").Value = 4562226.5
$ws.Range("K58").Value = 1358.8096
$ws.Range("L58").Value = 4562226.5
$ws.Range("M58").Value = -1155.8096
$ws.Range("N58").Value = -4562632.5
$ws.Range("H107").Value = 1249.2413
$ws.Range("I107").Value = 487.78946
$ws.Range("K107").Value = 487.78946
$ws.Range("M107").Value = 1432.21054
$ws.Range("H136").Value = 1569157.1
$ws.Range("I136").Value = 1358.8096
$ws.Range("J136").Value = 4562226.5
$ws.Range("K136").Value = 4076.4288
$ws.Range("L136").Value = 13686679.5
$ws.Range("M136").Value = -1526.4288
$ws.Range("N136").Value = -13691779.5

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 9293.637000000001
$ws.Range("I23").Value = 60
$ws.Range("J23").Value = 16988.334
$ws.Range("K23").Value = 180
$ws.Range("L23").Value = 50965.00199999999
$ws.Range("M23").Value = 55
$ws.Range("N23").Value = -51435.00199999999
$ws.Range("H68").Value = 576369.9
$ws.Range("I68").Value = 828.5
$ws.Range("J68").Value = 1299907.6
$ws.Range("K68").Value = 2485.5
$ws.Range("L68").Value = 3899722.8
$ws.Range("M68").Value = -1674.5
$ws.Range("N68").Value = -3901344.8
$ws.Range("H71").Value = 576369.9
$ws.Range("I71").Value = 828.5
$ws.Range("J71").Value = 1299907.6
$ws.Range("K71").Value = 7456.5
$ws.Range("L71").Value = 11699168.4
$ws.Range("M71").Value = -3400.5
$ws.Range("N71").Value = -11707280.4
$ws.Range("H107").Value = 13091547
$ws.Range("J107").Value = 458456.44
$ws.Range("L107").Value = 1375369.32
$ws.Range("N107").Value = -1379209.32
$ws.Range("H113").Value = 471.06897
$ws.Range("I113").Value = 441.8
$ws.Range("J113").Value = 502.42856
$ws.Range("K113").Value = 1325.4
$ws.Range("L113").Value = 1507.28568
$ws.Range("M113").Value = 844.5999999999999
$ws.Range("N113").Value = -5847.28568
$ws.Range("H132").Value = 1081
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 1081
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 9729
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -14789

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5222.5
$ws.Range("I70").Value = 5718
$ws.Range("K70").Value = 5718
$ws.Range("M70").Value = -5448
$ws.Range("H73").Value = 5222.5
$ws.Range("I73").Value = 5718
$ws.Range("K73").Value = 5718
$ws.Range("M73").Value = -4782
$ws.Range("H107").Value = 822.29034
$ws.Range("I107").Value = 466.86957
$ws.Range("J107").Value = 1844.125
$ws.Range("K107").Value = 466.86957
$ws.Range("L107").Value = 1844.125
$ws.Range("M107").Value = 1453.13043
$ws.Range("N107").Value = -5684.125

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 142861760
$ws.Range("I40").Value = 333334660
$ws.Range("J40").Value = 7076.25
$ws.Range("K40").Value = 333334660
$ws.Range("L40").Value = 7076.25
$ws.Range("M40").Value = -333334524
$ws.Range("N40").Value = -7348.25
$ws.Range("H46").Value = 910076.75
$ws.Range("I46").Value = 1007.4286
$ws.Range("J46").Value = 2500948
$ws.Range("K46").Value = 1007.4286
$ws.Range("L46").Value = 2500948
$ws.Range("M46").Value = -819.4286
$ws.Range("N46").Value = -2501324
$ws.Range("H122").Value = 14427.083
$ws.Range("I122").Value = 23080
$ws.Range("J122").Value = 8246.429
$ws.Range("K122").Value = 69240
$ws.Range("L122").Value = 24739.287
$ws.Range("M122").Value = -66790
$ws.Range("N122").Value = -29639.287
$ws.Range("H132").Value = 50007460
$ws.Range("I132").Value = 66673300
$ws.Range("J132").Value = 9940.6
$ws.Range("K132").Value = 200019900
$ws.Range("L132").Value = 29821.8
$ws.Range("M132").Value = -200017370
$ws.Range("N132").Value = -34881.8
$ws.Range("H136").Value = 20835490
$ws.Range("I136").Value = 71429970
$ws.Range("J136").Value = 2468.8235
$ws.Range("K136").Value = 214289910
$ws.Range("L136").Value = 7406.470499999999
$ws.Range("M136").Value = -214287360
$ws.Range("N136").Value = -12506.4705

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 8316.571
$ws.Range("I122").Value = 11910.667
$ws.Range("J122").Value = 3524.4443
$ws.Range("K122").Value = 35732.001
$ws.Range("L122").Value = 10573.3329
$ws.Range("M122").Value = -33282.001
$ws.Range("N122").Value = -15473.3329
$ws.Range("H132").Value = 2641.5557
$ws.Range("I132").Value = 2023.4706
$ws.Range("J132").Value = 3194.5789
$ws.Range("K132").Value = 6070.4118
$ws.Range("L132").Value = 9583.736699999999
$ws.Range("M132").Value = -3540.4118
$ws.Range("N132").Value = -14643.7367
$ws.Range("H136").Value = 12514203
$ws.Range("I136").Value = 22751416
$ws.Range("J136").Value = 2053.889
$ws.Range("K136").Value = 68254248
$ws.Range("L136").Value = 6161.667
$ws.Range("M136").Value = -68251698
$ws.Range("N136").Value = -11261.667
